$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2373.2307
$ws.Range("I6").Value = 51
$ws.Range("J6").Value = 2566.75
$ws.Range("K6").Value = 153
$ws.Range("L6").Value = 7700.25
$ws.Range("M6").Value = -41
$ws.Range("N6").Value = -7924.25
$ws.Range("H29").Value = 1500
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 629.5714
$ws.Range("I38").Value = 629.5714
$ws.Range("K38").Value = 1888.7142
$ws.Range("M38").Value = -1516.7142
$ws.Range("H39").Value = 925.0833
$ws.Range("I39").Value = 260.6
$ws.Range("J39").Value = 1399.7142
$ws.Range("K39").Value = 781.8000000000001
$ws.Range("L39").Value = 4199.142599999999
$ws.Range("M39").Value = -485.8000000000001
$ws.Range("N39").Value = -4791.142599999999
$ws.Range("H42").Value = 504.0909
$ws.Range("J42").Value = 650.875
$ws.Range("L42").Value = 1952.625
$ws.Range("N42").Value = -2412.625
$ws.Range("H76").Value = 5858481
$ws.Range("H79").Value = 5858481
$ws.Range("H132").Value = 908.44116
$ws.Range("I132").Value = 877.2239
$ws.Range("K132").Value = 2631.6717
$ws.Range("M132").Value = -101.6716999999999
$ws.Range("H137").Value = 975.675
$ws.Range("I137").Value = 696.2941
$ws.Range("K137").Value = 2088.8823
$ws.Range("M137").Value = 461.1177000000002
$ws.Range("H138").Value = 2311.2122
$ws.Range("I138").Value = 2402.5
$ws.Range("J138").Value = 1800
$ws.Range("K138").Value = 7207.5
$ws.Range("L138").Value = 5400
$ws.Range("M138").Value = -2067.5
$ws.Range("N138").Value = -15680
$ws.Range("H140").Value = 82608.95
$ws.Range("J140").Value = 82608.95
$ws.Range("L140").Value = 82608.95
$ws.Range("N140").Value = -92968.95

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5773.4863
$ws.Range("I32").Value = 4759.276
$ws.Range("J32").Value = 9450
$ws.Range("K32").Value = 4759.276
$ws.Range("L32").Value = 9450
$ws.Range("M32").Value = -4472.276
$ws.Range("N32").Value = -10024
$ws.Range("H61").Value = 2351.6333
$ws.Range("I61").Value = 1948.0769
$ws.Range("J61").Value = 4974.75
$ws.Range("K61").Value = 1948.0769
$ws.Range("L61").Value = 4974.75
$ws.Range("M61").Value = -1736.0769
$ws.Range("N61").Value = -5398.75
$ws.Range("H74").Value = 1205.9412
$ws.Range("I74").Value = 932.7308
$ws.Range("J74").Value = 2093.875
$ws.Range("K74").Value = 932.7308
$ws.Range("L74").Value = 2093.875
$ws.Range("M74").Value = -58.73080000000004
$ws.Range("N74").Value = -3841.875
$ws.Range("H77").Value = 1205.9412
$ws.Range("I77").Value = 932.7308
$ws.Range("J77").Value = 2093.875
$ws.Range("K77").Value = 4663.654
$ws.Range("L77").Value = 10469.375
$ws.Range("M77").Value = -295.6540000000005
$ws.Range("N77").Value = -19205.375
$ws.Range("H101").Value = 50683.715
$ws.Range("J101").Value = 50683.715
$ws.Range("L101").Value = 50683.715
$ws.Range("N101").Value = -57173.715
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
$ws.Range("H136").Value = 2351.6333
$ws.Range("I136").Value = 1948.0769
$ws.Range("J136").Value = 4974.75
$ws.Range("K136").Value = 5844.2307
$ws.Range("L136").Value = 14924.25
$ws.Range("M136").Value = -3294.2307
$ws.Range("N136").Value = -20024.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2800
$ws.Range("J99").Value = 2800
$ws.Range("L99").Value = 2800
$ws.Range("N99").Value = -5796
$ws.Range("H134").Value = 5154.357
$ws.Range("I134").Value = 6984.154
$ws.Range("K134").Value = 20952.462
$ws.Range("M134").Value = -18417.462

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1611994.1
$ws.Range("I58").Value = 2416873.8
$ws.Range("J58").Value = 2234.889
$ws.Range("K58").Value = 2416873.8
$ws.Range("L58").Value = 2234.889
$ws.Range("M58").Value = -2416670.8
$ws.Range("N58").Value = -2640.889
$ws.Range("H105").Value = 2338.5
$ws.Range("I105").Value = 2004
$ws.Range("J105").Value = 4011
$ws.Range("K105").Value = 2004
$ws.Range("L105").Value = 4011
$ws.Range("M105").Value = -257
$ws.Range("N105").Value = -7505
$ws.Range("H136").Value = 1611994.1
$ws.Range("I136").Value = 2416873.8
$ws.Range("J136").Value = 2234.889
$ws.Range("K136").Value = 7250621.399999999
$ws.Range("L136").Value = 6704.667
$ws.Range("M136").Value = -7248071.399999999
$ws.Range("N136").Value = -11804.667

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1088.3158
$ws.Range("J122").Value = 1144.6666
$ws.Range("L122").Value = 10301.9994
$ws.Range("N122").Value = -15201.9994
$ws.Range("H132").Value = 899.3125
$ws.Range("I132").Value = 692.6
$ws.Range("J132").Value = 1243.8334
$ws.Range("K132").Value = 6233.400000000001
$ws.Range("L132").Value = 11194.5006
$ws.Range("M132").Value = -3703.400000000001
$ws.Range("N132").Value = -16254.5006

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 21360
$ws.Range("J46").Value = 21360
$ws.Range("L46").Value = 21360
$ws.Range("N46").Value = -21672
$ws.Range("H80").Value = 3674.25
$ws.Range("I80").Value = 2199
$ws.Range("J80").Value = 4166
$ws.Range("K80").Value = 2199
$ws.Range("L80").Value = 4166
$ws.Range("M80").Value = -1201
$ws.Range("N80").Value = -6162
$ws.Range("H83").Value = 3674.25
$ws.Range("I83").Value = 2199
$ws.Range("J83").Value = 4166
$ws.Range("K83").Value = 10995
$ws.Range("L83").Value = 20830
$ws.Range("M83").Value = -6003
$ws.Range("N83").Value = -30814
$ws.Range("H132").Value = 803445.8
$ws.Range("I132").Value = 1375182.1
$ws.Range("J132").Value = 3014.95
$ws.Range("K132").Value = 4125546.3
$ws.Range("L132").Value = 9044.849999999999
$ws.Range("M132").Value = -4123016.3
$ws.Range("N132").Value = -14104.85

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1287.4615
$ws.Range("I82").Value = 1175.1
$ws.Range("J82").Value = 1662
$ws.Range("K82").Value = 1175.1
$ws.Range("L82").Value = 1662
$ws.Range("M82").Value = -814.0999999999999
$ws.Range("N82").Value = -2384
$ws.Range("H85").Value = 1287.4615
$ws.Range("I85").Value = 1175.1
$ws.Range("J85").Value = 1662
$ws.Range("K85").Value = 1175.1
$ws.Range("L85").Value = 1662
$ws.Range("M85").Value = 72.90000000000009
$ws.Range("N85").Value = -4158
$ws.Range("H136").Value = 2760.4688
$ws.Range("I136").Value = 1934.2632
$ws.Range("K136").Value = 5802.7896
$ws.Range("M136").Value = -3252.7896

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7998.8887
$ws.Range("I126").Value = 8298.571
$ws.Range("K126").Value = 24895.713
$ws.Range("M126").Value = -22425.713
$ws.Range("H136").Value = 13229147
$ws.Range("I136").Value = 16341383
$ws.Range("K136").Value = 49024149
$ws.Range("M136").Value = -49021599
